$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1
$ws.Range("H1").Value = "Save"

# Copy the formatting from the neighboring header cell (G1) onto H1 so the
# new header matches the existing bold/centered/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
